# Update gh-pages to output generated at 456a3b4
# Applies small numeric corrections to the "F" (follower/fan count) column
# across the "展览", "演出" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 9917
$ws1.Range("F5").Value = 693
$ws1.Range("F26").Value = 92
$ws1.Range("F32").Value = 3765
$ws1.Range("F33").Value = 781
$ws1.Range("F38").Value = 208

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 188

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 9917
$ws4.Range("F7").Value = 693
$ws4.Range("F26").Value = 92
$ws4.Range("F31").Value = 3765
$ws4.Range("F32").Value = 781
$ws4.Range("F37").Value = 208
